$wb = $excel.ActiveWorkbook

$newGuid = "2c251b13-2c98-4c63-b0ca-0238d039647f"
$newHash = "70a13ac1e3bd67aa4ab1f5a86e2fc9d7231af888"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

# Drop the existing hyperlink(s) on this sheet so we can safely rewrite B2's
# display text (this runtime's Hyperlink.TextToDisplay setter duplicates the
# entry rather than updating it in place, so delete+recreate is used instead).
$wsOverview.Range("B2").Hyperlinks.Delete()

$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("G2").Value = "2016-11-15 17:35:12"

$ovB2 = $wsOverview.Range("B2")
$wsOverview.Hyperlinks.Add($ovB2, "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec1c347d7faa238ebde65be360d8e7bf68b39427/e2e/$newGuid.md", "", "", "e2e\$newGuid.md") | Out-Null

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Remove all hyperlinks on the sheet (A2 + I2) so the I2 one can be dropped;
# A2's is recreated below with the refreshed display text.
$wsZh.Range("I2").Hyperlinks.Delete()

$wsZh.Range("A2").Value = "$newGuid.md"
$wsZh.Range("G2").Value = "$newGuid.$newHash.zh-cn.xlf"
$wsZh.Range("H2").Value = "2016-11-15 17:34:59"
$wsZh.Range("K2").Value = "0001-01-01 00:00:00"

# Clear "Latest Target File" (I2) and "Latest Handback File" (J2)
$zhI2 = $wsZh.Range("I2")
$zhI2.Value = ""
$zhI2.Style = "Normal"
$wsZh.Range("J2").Value = ""

# "Has metadata": False -> True. Assigning the literal string "True" via
# .Value gets auto-coerced to a native boolean cell by this runtime, so copy
# it from an existing text cell that already holds "True" to keep it a
# shared-string text cell (matching the source workbook's modelling).
$wsZh.Range("M2").Copy($wsZh.Range("O2"))

$zhA2 = $wsZh.Range("A2")
$wsZh.Hyperlinks.Add($zhA2, "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/ec1c347d7faa238ebde65be360d8e7bf68b39427/e2e/$newGuid.md", "", "", "$newGuid.md") | Out-Null

# Column widths for "Latest Target File" / "Latest Handback File"
$wsZh.Columns.Item(9).ColumnWidth = 17.8
$wsZh.Columns.Item(10).ColumnWidth = 20.8

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("I2").Hyperlinks.Delete()

$wsDe.Range("A2").Value = "$newGuid.md"
$wsDe.Range("G2").Value = "$newGuid.$newHash.de-de.xlf"
$wsDe.Range("H2").Value = "2016-11-15 17:35:12"
$wsDe.Range("K2").Value = "0001-01-01 00:00:00"

$deI2 = $wsDe.Range("I2")
$deI2.Value = ""
$deI2.Style = "Normal"
$wsDe.Range("J2").Value = ""

$wsDe.Range("M2").Copy($wsDe.Range("O2"))

$deA2 = $wsDe.Range("A2")
$wsDe.Hyperlinks.Add($deA2, "https://github.com/OpenLocalizationTestOrg/ol-test0-dede/blob/961c924b108f55c424dcc16385e8aba5f1224256/e2e/$newGuid.md", "", "", "$newGuid.md") | Out-Null

$wsDe.Columns.Item(9).ColumnWidth = 17.8
$wsDe.Columns.Item(10).ColumnWidth = 20.8
